$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.098.69"
$ws.Range("E2").Value = "  -0.57%  "
$ws.Range("D3").Value = "1.893.65"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.46"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5214"
$ws.Range("E7").Value = "  -0.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3756"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07261"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.12"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8976"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08187"
$ws.Range("E12").Value = "  +0.31%  "
$ws.Range("D13").Value = "1.945.20"
$ws.Range("E13").Value = "  +2.03%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "96.26"
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.295"
$ws.Range("E15").Value = "  +0.00%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008587"
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.55"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "27.134.98"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.69"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.404"
$ws.Range("E23").Value = "  -0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "148.58"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.18"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "115.09"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.787"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.858"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09209"
$ws.Range("E31").Value = "  -0.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05032"
$ws.Range("E32").Value = "  -0.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7885"
$ws.Range("E33").Value = "  -2.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.214"
$ws.Range("E34").Value = "  -2.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.424"
$ws.Range("E35").Value = "  +2.81%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.971"
$ws.Range("E36").Value = "  -0.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.604"
$ws.Range("E37").Value = "  +0.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5712"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01990"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("E40").Value = "  -0.50%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.998"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.543"
$ws.Range("E42").Value = "  -1.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.35"
$ws.Range("E43").Value = "  -2.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1511"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4845"
$ws.Range("E45").Value = "  -0.39%  "
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.03"
$ws.Range("E47").Value = "  -1.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.617"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.12"
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "63.52"
$ws.Range("E51").Value = "  -0.33%  "
